# Update scripts with new TPM-based NATMI values.
# The underlying analysis was re-run; the "ECs" target cluster no longer
# appears, and the remaining Sending-cluster / Resolving-Mac rows were
# recomputed with new statistics. Rows that used to target "ECs" are
# removed, and the surviving rows (originally rows 3, 5 and 7, all
# targeting "Resolving-Mac") move up to rows 2-4 with refreshed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows that are no longer present (old rows 5, 6 and 7).
# Delete from the bottom up so row numbers of rows still to be removed
# don't shift.
$ws.Rows("7").Delete()
$ws.Rows("6").Delete()
$ws.Rows("5").Delete()

function Set-RowValues($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# Row 2: FAPs -> Cxcl13 -> Cxcr3 -> Resolving-Mac
Set-RowValues 2 @(
    "FAPs", "Cxcl13", "Cxcr3", "Resolving-Mac",
    2, 0.6666666666666666, 0.2412233333333333, 0.72367,
    0.2038054651530871, 0.2038054651530872,
    3, 1, 1.888791333333333, 5.666374, 1, 1,
    0.4556205413977779, 4.100584872580001,
    0.2038054651530871, 0.2038054651530872
)

# Row 3: MuSCs -> Cxcl13 -> Cxcr3 -> Resolving-Mac
Set-RowValues 3 @(
    "MuSCs", "Cxcl13", "Cxcr3", "Resolving-Mac",
    1, 0.3333333333333333, 0.133983, 0.401949,
    0.1131999432238703, 0.1131999432238703,
    3, 1, 1.888791333333333, 5.666374, 1, 1,
    0.253065929214, 2.277593362926,
    0.1131999432238703, 0.1131999432238703
)

# Row 4: Resolving-Mac -> Cxcl13 -> Cxcr3 -> Resolving-Mac
Set-RowValues 4 @(
    "Resolving-Mac", "Cxcl13", "Cxcr3", "Resolving-Mac",
    3, 1, 0.8083896666666667, 2.425169,
    0.6829945916230425, 0.6829945916230425,
    3, 1, 1.888791333333333, 5.666374, 1, 1,
    1.526879396356222, 13.741914567206,
    0.6829945916230425, 0.6829945916230425
)
